$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $value) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws "D2" '28.606.11'
Set-TextValue $ws "E2" '  +1.39%  '
Set-TextValue $ws "D3" '1.793.71'
Set-TextValue $ws "E3" '  -0.41%  '
Set-TextValue $ws "D4" '1.003'
Set-TextValue $ws "E4" '  -0.08%  '
Set-TextValue $ws "D5" '329.34'
Set-TextValue $ws "E5" '  -2.72%  '
Set-TextValue $ws "E6" '  +0.12%  '
Set-TextValue $ws "D7" '0.4406'
Set-TextValue $ws "E7" '  -2.92%  '
Set-TextValue $ws "D8" '0.3742'
Set-TextValue $ws "E8" '  +5.66%  '
Set-TextValue $ws "D9" '45.70'
Set-TextValue $ws "E9" '  +0.18%  '
Set-TextValue $ws "D10" '0.07610'
Set-TextValue $ws "E10" '  +0.53%  '
Set-TextValue $ws "D11" '1.135'
Set-TextValue $ws "E11" '  -1.37%  '
Set-TextValue $ws "D12" '22.69'
Set-TextValue $ws "E12" '  -0.39%  '
Set-TextValue $ws "D13" '1.005'
Set-TextValue $ws "E13" '  +0.11%  '
Set-TextValue $ws "D14" '6.233'
Set-TextValue $ws "E14" '  -0.34%  '
Set-TextValue $ws "D15" '7.499'
Set-TextValue $ws "E15" '  +3.06%  '
Set-TextValue $ws "D16" '1.793.02'
Set-TextValue $ws "E16" '  -0.37%  '
Set-TextValue $ws "D17" '0.00001089'
Set-TextValue $ws "E17" '  -0.15%  '
Set-TextValue $ws "D18" '0.06701'
Set-TextValue $ws "E18" '  +0.14%  '
Set-TextValue $ws "D19" '80.65'
Set-TextValue $ws "E20" '  +0.23%  '
Set-TextValue $ws "E21" '  +1.27%  '
Set-TextValue $ws "D22" '6.196'
Set-TextValue $ws "E22" '  -3.68%  '
Set-TextValue $ws "D23" '28.560.28'
Set-TextValue $ws "E23" '  +1.28%  '
Set-TextValue $ws "D24" '11.73'
Set-TextValue $ws "E24" '  -2.23%  '
Set-TextValue $ws "D25" '2.442'
Set-TextValue $ws "E25" '  +2.14%  '
Set-TextValue $ws "D26" '20.43'
Set-TextValue $ws "E26" '  -1.63%  '
Set-TextValue $ws "D27" '153.31'
Set-TextValue $ws "E27" '  -1.63%  '
Set-TextValue $ws "D28" '2.342'
Set-TextValue $ws "E28" '  -3.36%  '
Set-TextValue $ws "D29" '1.999.08'
Set-TextValue $ws "E29" '  -0.29%  '
Set-TextValue $ws "D30" '1.314'
Set-TextValue $ws "E30" '  +1.61%  '
Set-TextValue $ws "D31" '130.94'
Set-TextValue $ws "E31" '  -1.97%  '
Set-TextValue $ws "E32" '  -2.37%  '
Set-TextValue $ws "D33" '5.804'
Set-TextValue $ws "E33" '  -2.38%  '
Set-TextValue $ws "D34" '0.09263'
Set-TextValue $ws "E34" '  -2.32%  '
Set-TextValue $ws "D35" '0.2254'
Set-TextValue $ws "E35" '  +3.41%  '
Set-TextValue $ws "D36" '12.13'
Set-TextValue $ws "E36" '  -0.52%  '
Set-TextValue $ws "D37" '0.06259'
Set-TextValue $ws "E37" '  -0.61%  '
Set-TextValue $ws "D38" '0.02318'
Set-TextValue $ws "E38" '  -2.89%  '
Set-TextValue $ws "D39" '5.208'
Set-TextValue $ws "D40" '0.6576'
Set-TextValue $ws "E40" '  -2.48%  '
Set-TextValue $ws "D41" '1.195'
Set-TextValue $ws "E41" '  -1.82%  '
Set-TextValue $ws "B42" 'FraxShare'
Set-TextValue $ws "C42" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws "D42" '7.986'
Set-TextValue $ws "E42" '  -2.15%  '
Set-TextValue $ws "B43" 'WEMIXTOKEN'
Set-TextValue $ws "C43" 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws "D43" '1.423'
Set-TextValue $ws "E43" '  -4.12%  '
Set-TextValue $ws "D45" '13.86'
Set-TextValue $ws "E45" '  -0.79%  '
Set-TextValue $ws "D46" '0.6088'
Set-TextValue $ws "E46" '  -0.67%  '
Set-TextValue $ws "D47" '3.810'
Set-TextValue $ws "E47" '  -1.54%  '
Set-TextValue $ws "D48" '127.57'
Set-TextValue $ws "E48" '  -1.67%  '
Set-TextValue $ws "D49" '2.016'
Set-TextValue $ws "E49" '  -1.18%  '
Set-TextValue $ws "D50" '0.07012'
Set-TextValue $ws "E50" '  -1.37%  '
Set-TextValue $ws "D51" '1.140'
Set-TextValue $ws "E51" '  -2.74%  '

"done"
